$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header formatting (bold, centered, thin border) by copying
# the format from the neighboring header cell H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# I0 / IF data values for rows 2-88
$iVals = @(6,6,7,9,7,7,8,7,7,9,8,10,7,8,8,8,8,8,8,7,6,8,8,8,8,8,7,9,8,7,7,8,8,7,7,8,9,7,8,8,8,8,8,8,8,9,7,8,8,7,7,7,9,7,7,7,7,8,8,8,9,8,7,8,7,9,8,6,9,7,8,7,8,8,7,6,8,8,8,7,6,5,8,8,5,5,5)
$jVals = @(6,7,7,9,7,8,8,8,7,9,8,10,8,8,8,8,8,8,8,8,7,8,9,8,8,8,7,9,8,8,8,8,8,8,7,8,9,8,8,8,8,8,8,8,8,10,7,8,8,8,8,7,9,7,7,8,7,9,8,8,9,8,7,8,7,9,8,6,9,8,8,8,8,9,7,7,8,8,9,7,6,6,8,8,5,5,5)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
